$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.03024452979876618
$ws.Range("H2").Value = -28.463838563026
$ws.Range("I2").Value = 101.733285878834
$ws.Range("G3").Value = 0.07104583732611526
$ws.Range("H3").Value = 45.88363016775652
$ws.Range("G4").Value = 0.01756330217041502
$ws.Range("H4").Value = 822.4929451059314
$ws.Range("G5").Value = 0.02620898725219534
$ws.Range("H5").Value = 590.3731447692578
$ws.Range("G6").Value = 0.0376453636828554
$ws.Range("H6").Value = 8.583271769210267
$ws.Range("G7").Value = 0.03745091098883087
$ws.Range("H7").Value = -29.58892121481232
$ws.Range("G8").Value = -0.00594064439204748
$ws.Range("H8").Value = 68.43570796497262
$ws.Range("G9").Value = -0.001385632454148331
$ws.Range("H9").Value = 93.56355180295316
$ws.Range("G10").Value = -0.06552774159793567
$ws.Range("H10").Value = 9.867236043281249
$ws.Range("G11").Value = -0.07587304609578838
$ws.Range("H11").Value = 17.53919439010665
$ws.Range("G12").Value = -0.2222606986682292
$ws.Range("H12").Value = 9.078962536213639
$ws.Range("G13").Value = -0.2974994537772763
$ws.Range("H13").Value = -8.255146086400647
$ws.Range("G14").Value = -0.05373465913273377
$ws.Range("H14").Value = -44.84831305619134
$ws.Range("G15").Value = 0.05431094315511857
$ws.Range("H15").Value = 256.194477982264
$ws.Range("G16").Value = 0.1235641973023537
$ws.Range("H16").Value = -1.391760732183239
$ws.Range("G17").Value = 0.1564427157106385
$ws.Range("H17").Value = 11.54521769741305
$ws.Range("G18").Value = 0.1194078031445488
$ws.Range("H18").Value = -4.269064837120395
$ws.Range("G19").Value = 0.1274193101274259
$ws.Range("H19").Value = -4.349797757830496
$ws.Range("G20").Value = 0.04021232460564812
$ws.Range("H20").Value = 17.11379485391262
$ws.Range("G21").Value = 0.04867663486258211
$ws.Range("H21").Value = -16.13240853702199
$ws.Range("G22").Value = -0.07158136933928702
$ws.Range("H22").Value = 10.35294111122793
$ws.Range("G23").Value = -0.09697757530516884
$ws.Range("H23").Value = -55.04892678921694
$ws.Range("G24").Value = 0.1177125661533811
$ws.Range("H24").Value = -0.3403001748906242
$ws.Range("G25").Value = 0.1118525185123732
$ws.Range("H25").Value = -11.34804672029393
$ws.Range("G26").Value = 0.0591177842344799
$ws.Range("H26").Value = 18.93943251680953
$ws.Range("G27").Value = 0.05948649479817009
$ws.Range("H27").Value = -31.36957935717876
$ws.Range("G28").Value = -0.08595534333522964
$ws.Range("H28").Value = -35.16258455085313
$ws.Range("G29").Value = -0.08241589393561426
$ws.Range("H29").Value = -15.7978388662375
$ws.Range("G30").Value = 0.0729807116861432
$ws.Range("H30").Value = 14.55479129374364
$ws.Range("G31").Value = 0.04535069509718408
$ws.Range("H31").Value = -25.13973292629875
$ws.Range("G32").Value = 0.06549111880355671
$ws.Range("H32").Value = -33.35185897615198
$ws.Range("G33").Value = 0.1017666693602952
$ws.Range("H33").Value = 23.67628855808154
$ws.Range("G34").Value = 0.014630848159418
$ws.Range("H34").Value = -43.84675920969583
$ws.Range("G35").Value = -0.01186505017964416
$ws.Range("H35").Value = -5.842416858380166
$ws.Range("G36").Value = 0.03007868271686307
$ws.Range("H36").Value = 5488.682081029435
$ws.Range("G37").Value = -0.003502189967209183
$ws.Range("H37").Value = 72.10358405226897
$ws.Range("G38").Value = 0.07278370639071388
$ws.Range("H38").Value = -32.14100253217061
$ws.Range("G39").Value = 0.09728538022048455
$ws.Range("H39").Value = 13.56785057077581
$ws.Range("G40").Value = 0.01190504196651193
$ws.Range("H40").Value = 300.8084422860321
$ws.Range("G41").Value = 0.0196958696081565
$ws.Range("H41").Value = 31.33428438199767
$ws.Range("G42").Value = 0.1235658145244171
$ws.Range("H42").Value = 22.42102427734318
$ws.Range("G43").Value = 0.1178434562667192
$ws.Range("H43").Value = -1.915232815270834
$ws.Range("G44").Value = 0.02383806093267719
$ws.Range("H44").Value = -33.2034948656252
$ws.Range("G45").Value = 0.03192040387194704
$ws.Range("H45").Value = 94.9912294962332
$ws.Range("G46").Value = 0.06466868648872473
$ws.Range("H46").Value = 78.45289492399084
$ws.Range("G47").Value = 0.05168877455710817
$ws.Range("H47").Value = 2.475124746494537
$ws.Range("G48").Value = 0.032104565856088
$ws.Range("H48").Value = -24.95654833658694
$ws.Range("G49").Value = 0.06167364336885237
$ws.Range("H49").Value = -11.23041913228128
$ws.Range("G50").Value = 0.0190815693623699
$ws.Range("H50").Value = 10.47237661358225
$ws.Range("G51").Value = 0.02507954258999711
$ws.Range("H51").Value = 28.81312172622845
$ws.Range("G52").Value = -0.1063797840509629
$ws.Range("H52").Value = -2.762090475252428
$ws.Range("G53").Value = -0.09509200897074017
$ws.Range("H53").Value = -2.963565021964672
$ws.Range("G54").Value = 0.071011525489643
$ws.Range("H54").Value = -2.887115318753921
$ws.Range("G55").Value = 0.07671074283154139
$ws.Range("H55").Value = 23.82355231897832
$ws.Range("G56").Value = 0.0341787461613418
$ws.Range("H56").Value = -2.31633883728726
$ws.Range("G57").Value = 0.0149703396837001
$ws.Range("H57").Value = 159.2932141657255
$ws.Range("G58").Value = 0.04302553525589588
$ws.Range("H58").Value = 72.02915284223023
$ws.Range("G59").Value = 0.03176422991046859
$ws.Range("H59").Value = 34.1472693768815
$ws.Range("G60").Value = 0.04135616179583994
$ws.Range("H60").Value = 27.47464094653223
$ws.Range("G61").Value = 0.04687787486306264
$ws.Range("H61").Value = 270.3420774066287
$ws.Range("G62").Value = 0.0610339271685276
$ws.Range("H62").Value = 1.112127474737234
$ws.Range("G63").Value = 0.06703571438024759
$ws.Range("H63").Value = 105.6973199035065
$ws.Range("G64").Value = 0.04097478374210759
$ws.Range("H64").Value = 1.106945487715235
$ws.Range("G65").Value = 0.04985649981749831
$ws.Range("H65").Value = -11.06887028526866
$ws.Range("G66").Value = 0.07866416971532111
$ws.Range("H66").Value = -15.91585418848062
$ws.Range("G67").Value = 0.1288651572240292
$ws.Range("H67").Value = 11.62311773384005
$ws.Range("G68").Value = -0.03808273750484675
$ws.Range("H68").Value = -9.27463371999886
$ws.Range("G69").Value = -0.02221914123697266
$ws.Range("H69").Value = -4.699346362983955
$ws.Range("G70").Value = 0.06236038185475334
$ws.Range("H70").Value = -32.68253356744713
$ws.Range("G71").Value = 0.08895172043278908
$ws.Range("H71").Value = -2.474625207928791
$ws.Range("G72").Value = -0.05460853430156384
$ws.Range("H72").Value = 2.625650671821967
$ws.Range("G73").Value = -0.06550506458061546
$ws.Range("H73").Value = 11.19473340953218
$ws.Range("G74").Value = 0.1098499763538468
$ws.Range("H74").Value = 9.908440175637542
$ws.Range("G75").Value = 0.1406386173343374
$ws.Range("H75").Value = 44.38686765763468
$ws.Range("G76").Value = -0.00352592518455871
$ws.Range("H76").Value = -113.7889819178554
$ws.Range("G77").Value = 0.01119260878414328
$ws.Range("H77").Value = -20.67065748936501
$ws.Range("G78").Value = 0.09364663940164723
$ws.Range("H78").Value = 45.69237071475042
$ws.Range("G79").Value = 0.09666924155290466
$ws.Range("H79").Value = 26.01268647758129
$ws.Range("G80").Value = -0.1675174004908132
$ws.Range("H80").Value = -1.152718178376996
$ws.Range("G81").Value = -0.1516918080082811
$ws.Range("H81").Value = 27.79744229682048
$ws.Range("G82").Value = 0.119646037355407
$ws.Range("H82").Value = 4.316400172078088
$ws.Range("G83").Value = 0.1764312836270396
$ws.Range("H83").Value = -0.8713954377009145
$ws.Range("G84").Value = 0.07095170537769502
$ws.Range("H84").Value = 197.6455658282839
$ws.Range("G85").Value = 0.07916975958736094
$ws.Range("H85").Value = 28.57255499800427
